$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) cells we touch to remain Text, matching the
# workbook author convention of storing these as literal strings even when
# they look numeric (so Excel does not silently convert "1.014" -> 1.014 number).
$ws.Range("D2,D3,D4,D5,D6,D7,D8,D9,D10,D12,D13,D14,D15,D16,D17,D18,D19,D20,D21,D22,D23,D24,D25,D26,D27,D28,D29,D30,D31,D32,D33,D34,D35,D36,D37,D38,D39,D40,D41,D42,D43,D44,D45,D46,D47,D48,D49,D50,D51").NumberFormat = "@"

$ws.Range("D2").Value = '28.461.76'
$ws.Range("E2").Value = '  -2.24%  '

$ws.Range("D3").Value = '1.950.88'
$ws.Range("E3").Value = '  -0.88%  '

$ws.Range("D4").Value = '1.014'
$ws.Range("E4").Value = '  +0.74%  '

$ws.Range("D5").Value = '322.10'
$ws.Range("E5").Value = '  -1.69%  '

$ws.Range("D6").Value = '1.012'
$ws.Range("E6").Value = '  +0.67%  '

$ws.Range("D7").Value = '0.4790'
$ws.Range("E7").Value = '  -4.26%  '

$ws.Range("D8").Value = '0.4083'
$ws.Range("E8").Value = '  -3.10%  '

$ws.Range("D9").Value = '53.96'
$ws.Range("E9").Value = '  +2.63%  '

$ws.Range("D10").Value = '0.08470'
$ws.Range("E10").Value = '  -8.12%  '

$ws.Range("E11").Value = '  -3.76%  '

$ws.Range("D12").Value = '22.41'
$ws.Range("E12").Value = '  -2.43%  '

$ws.Range("D13").Value = '1.979.94'
$ws.Range("E13").Value = '  -0.21%  '

$ws.Range("D14").Value = '7.582'
$ws.Range("E14").Value = '  -3.76%  '

$ws.Range("D15").Value = '6.166'
$ws.Range("E15").Value = '  -4.41%  '

$ws.Range("D16").Value = '1.015'
$ws.Range("E16").Value = '  +0.87%  '

$ws.Range("D17").Value = '90.55'
$ws.Range("E17").Value = '  -1.15%  '

$ws.Range("D18").Value = '0.00001074'
$ws.Range("E18").Value = '  -2.70%  '

$ws.Range("D19").Value = '0.06638'
$ws.Range("E19").Value = '  -1.24%  '

$ws.Range("D20").Value = '18.44'
$ws.Range("E20").Value = '  -4.30%  '

$ws.Range("D21").Value = '1.012'
$ws.Range("E21").Value = '  +0.71%  '

$ws.Range("D22").Value = '5.841'
$ws.Range("E22").Value = '  -2.21%  '

$ws.Range("D23").Value = '28.470.30'
$ws.Range("E23").Value = '  -2.28%  '

$ws.Range("D24").Value = '11.45'
$ws.Range("E24").Value = '  -5.77%  '

$ws.Range("D25").Value = '2.303'
$ws.Range("E25").Value = '  +0.81%  '

$ws.Range("D26").Value = '2.173.39'
$ws.Range("E26").Value = '  -1.93%  '

$ws.Range("D27").Value = '156.43'
$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("D28").Value = '20.32'
$ws.Range("E28").Value = '  -1.61%  '

$ws.Range("D29").Value = '2.171'
$ws.Range("E29").Value = '  -4.25%  '

$ws.Range("D30").Value = '5.817'
$ws.Range("E30").Value = '  -6.05%  '

$ws.Range("D31").Value = '124.44'
$ws.Range("E31").Value = '  -1.88%  '

$ws.Range("D32").Value = '0.9857'
$ws.Range("E32").Value = '  -5.61%  '

$ws.Range("D33").Value = '0.09669'
$ws.Range("E33").Value = '  -1.84%  '

$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '3.697'
$ws.Range("E34").Value = '  +0.64%  '

$ws.Range("B35").Value = 'Filecoin'
$ws.Range("C35").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D35").Value = '5.636'
$ws.Range("E35").Value = '  -2.54%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '1.442'
$ws.Range("E36").Value = '  -5.92%  '

$ws.Range("D37").Value = '9.116'
$ws.Range("E37").Value = '  +1.67%  '

$ws.Range("D38").Value = '0.02334'
$ws.Range("E38").Value = '  -3.62%  '

$ws.Range("D39").Value = '0.06203'
$ws.Range("E39").Value = '  -2.28%  '

$ws.Range("D40").Value = '1.247'
$ws.Range("E40").Value = '  -4.09%  '

$ws.Range("D41").Value = '0.6225'
$ws.Range("E41").Value = '  -3.56%  '

$ws.Range("D42").Value = '11.21'
$ws.Range("E42").Value = '  -2.16%  '

$ws.Range("D43").Value = '1.012'
$ws.Range("E43").Value = '  +0.77%  '

$ws.Range("D44").Value = '0.1920'
$ws.Range("E44").Value = '  -3.33%  '

$ws.Range("D45").Value = '1.338'
$ws.Range("E45").Value = '  +4.51%  '

$ws.Range("D46").Value = '0.5954'
$ws.Range("E46").Value = '  -4.61%  '

$ws.Range("D47").Value = '12.98'
$ws.Range("E47").Value = '  -3.27%  '

$ws.Range("D48").Value = '2.063'
$ws.Range("E48").Value = '  -5.85%  '

$ws.Range("D49").Value = '3.413'
$ws.Range("E49").Value = '  -1.69%  '

$ws.Range("D50").Value = '0.06829'
$ws.Range("E50").Value = '  -2.01%  '

$ws.Range("D51").Value = '110.85'
$ws.Range("E51").Value = '  -1.96%  '
